$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 361, shifting rows 361:375 down to 362:376
$ws.Rows.Item(361).Insert(-4121)

# Populate the newly inserted row 361 with the new weekly record
$ws.Cells.Item(361, 1).Value = 6
$ws.Cells.Item(361, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(361, 3).Value = "Metropolitana"
$ws.Cells.Item(361, 4).Value = 44509
$ws.Cells.Item(361, 4).NumberFormat = $ws.Cells.Item(362, 4).NumberFormat
$ws.Cells.Item(361, 5).Value = 13
$ws.Cells.Item(361, 6).Value = 100112044
$ws.Cells.Item(361, 7).Value = "Perejil"
$ws.Cells.Item(361, 8).Value = "Sin especificar"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 240
$ws.Cells.Item(361, 11).Value = 11000
$ws.Cells.Item(361, 12).Value = 12000
$ws.Cells.Item(361, 13).Value = 11458
$ws.Cells.Item(361, 14).Value = "$/docena de atados"
$ws.Cells.Item(361, 15).Value = "Región Metropolitana"
$ws.Cells.Item(361, 16).Value = 3819
$ws.Cells.Item(361, 17).Value = 3
$ws.Cells.Item(361, 18).Value = "Hortaliza"
